$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$ws.Cells.Item(5, 3).Value = "FAIL"
$ws.Cells.Item(5, 4).Value = "No exception thrown"

$ws.Cells.Item(19, 3).Value = "FAIL"
$ws.Cells.Item(19, 4).Value = "No exception thrown"

$ws.Cells.Item(20, 3).Value = "FAIL"
$ws.Cells.Item(20, 4).Value = "No exception thrown"
